$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Task "4. Include 'Link Generator' as a feature" (row 5):
#   Status: Open -> In progress
#   Assignee: (blank) -> Arthur
$ws.Range("C5").Value = "In progress"
$ws.Range("D5").Value = "Arthur"

# Task "5. Search algorithm optimization" (row 6):
#   Status: Reopened to following optimization -> Done
#   (reuse the same "Done" cell formatting already used elsewhere, e.g. C4)
$ws.Range("C6").Value = "Done"
$ws.Range("C4").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the cursor where the author last left it
$ws.Range("C15").Select()
